# Restore C10 ("R30" rule's lower bound) value from 18 back to 1,
# as part of restoring the sheet to the prior committed revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
